$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.759.88"
$ws.Range("E2").Value = "  +2.08%  "
$ws.Range("D3").Value = "1.857.89"
$ws.Range("E3").Value = "  +1.78%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("E5").Value = "  +1.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6453"
$ws.Range("E6").Value = "  +3.92%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07510"
$ws.Range("E8").Value = "  +2.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2995"
$ws.Range("E9").Value = "  +3.38%  "
$ws.Range("B10").Value = "Solana"
$ws.Range("C10").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.13"
$ws.Range("E10").Value = "  +4.88%  "
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07692"
$ws.Range("E11").Value = "  +0.21%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.872.19"
$ws.Range("E12").Value = "  +2.62%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.055"
$ws.Range("E13").Value = "  +1.93%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6846"
$ws.Range("E14").Value = "  +2.97%  "
$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "84.00"
$ws.Range("E15").Value = "  +2.03%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009440"
$ws.Range("E16").Value = "  +5.05%  "
$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.068"
$ws.Range("E17").Value = "  +3.79%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "29.738.45"
$ws.Range("E18").Value = "  +2.09%  "
$ws.Range("B19").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C19").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D19").Value = "2.122.58"
$ws.Range("E19").Value = "  +2.49%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "240.29"
$ws.Range("E20").Value = "  +0.60%  "
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.69"
$ws.Range("E21").Value = "  +2.12%  "
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.451"
$ws.Range("E23").Value = "  +2.84%  "
$ws.Range("B24").Value = "BinanceUSD"
$ws.Range("C24").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.002"
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.85"
$ws.Range("E25").Value = "  +0.53%  "
$ws.Range("B26").Value = "Stellar"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1430"
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.543"
$ws.Range("E27").Value = "  +0.70%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.97"
$ws.Range("E28").Value = "  +1.69%  "
$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06157"
$ws.Range("E29").Value = "  +10.50%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.504"
$ws.Range("E30").Value = "  +1.28%  "
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.270"
$ws.Range("E31").Value = "  +5.48%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.154"
$ws.Range("E32").Value = "  +1.45%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.129"
$ws.Range("E33").Value = "  +0.96%  "
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.886"
$ws.Range("E34").Value = "  +2.40%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.158"
$ws.Range("E35").Value = "  +2.10%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7348"
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.609"
$ws.Range("E37").Value = "  -0.68%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.855"
$ws.Range("E38").Value = "  +0.59%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01797"
$ws.Range("E39").Value = "  +1.86%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "1.217.72"
$ws.Range("E40").Value = "  +0.35%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9304"
$ws.Range("E41").Value = "  +1.97%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.221"
$ws.Range("E42").Value = "  -1.31%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.002"
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "2.032.50"
$ws.Range("E44").Value = "  +2.90%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "102.05"
$ws.Range("E45").Value = "  +0.40%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "66.45"
$ws.Range("E46").Value = "  +2.77%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000122"
$ws.Range("E47").Value = "  +3.84%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5066"
$ws.Range("E48").Value = "  -0.38%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.338"
$ws.Range("E49").Value = "  +2.39%  "
$ws.Range("B50").Value = "TheSandbox"
$ws.Range("C50").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4085"
$ws.Range("E50").Value = "  +1.60%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1134"
$ws.Range("E51").Value = "  +1.79%  "
